$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/160fbf46670014aeebc0e0075d4d72572dd46f63/e2e/b72a1f6e-20e8-451a-8868-4fc3456e35e5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f00c2a7464b76d295e63e58c3e666bb95b4aa119/e2e/b72a1f6e-20e8-451a-8868-4fc3456e35e5.md."
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f00c2a7464b76d295e63e58c3e666bb95b4aa119/e2e/b72a1f6e-20e8-451a-8868-4fc3456e35e5.md"
$displayName = "b72a1f6e-20e8-451a-8868-4fc3456e35e5.md"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "b72a1f6e-20e8-451a-8868-4fc3456e35e5.a5be7e43187da8aa0a5d861bf03b56f08af34128.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-06 06:50:24"
$wsZh.Range("P7").Value = $errorMessage
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestUrl, $null, $null, $displayName)

# widen the "Error Detail" column (column 16 / P) to fit the new long message
$wsZh.Columns.Item(16).ColumnWidth = $wsZh.Columns.Item(1).ColumnWidth

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "b72a1f6e-20e8-451a-8868-4fc3456e35e5.a5be7e43187da8aa0a5d861bf03b56f08af34128.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-06 06:50:32"
$wsDe.Range("P7").Value = $errorMessage
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestUrl, $null, $null, $displayName)

# widen the "Error Detail" column (column 16 / P) to fit the new long message
$wsDe.Columns.Item(16).ColumnWidth = $wsDe.Columns.Item(1).ColumnWidth
